$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1688.8
$ws.Range("I28").Value = 1572.9412
$ws.Range("J28").Value = 1840.3077
$ws.Range("K28").Value = 1572.9412
$ws.Range("L28").Value = 1840.3077
$ws.Range("M28").Value = -1087.9412
$ws.Range("N28").Value = -2810.3077
$ws.Range("H62").Value = 3426.5908
$ws.Range("I62").Value = 1956.6428
$ws.Range("J62").Value = 5999
$ws.Range("K62").Value = 1956.6428
$ws.Range("L62").Value = 5999
$ws.Range("M62").Value = -1332.6428
$ws.Range("N62").Value = -7247
$ws.Range("H65").Value = 3426.5908
$ws.Range("I65").Value = 1956.6428
$ws.Range("J65").Value = 5999
$ws.Range("K65").Value = 9783.214
$ws.Range("L65").Value = 29995
$ws.Range("M65").Value = -6663.214
$ws.Range("N65").Value = -36235
$ws.Range("H76").Value = 3485.7144
$ws.Range("I76").Value = 3650
$ws.Range("J76").Value = 3266.6667
$ws.Range("K76").Value = 3650
$ws.Range("L76").Value = 3266.6667
$ws.Range("M76").Value = -3335
$ws.Range("N76").Value = -3896.6667
$ws.Range("H79").Value = 3485.7144
$ws.Range("I79").Value = 3650
$ws.Range("J79").Value = 3266.6667
$ws.Range("K79").Value = 3650
$ws.Range("L79").Value = 3266.6667
$ws.Range("M79").Value = -2558
$ws.Range("N79").Value = -5450.6667
$ws.Range("H86").Value = 4084.5
$ws.Range("I86").Value = 3160.6667
$ws.Range("J86").Value = 4840.364
$ws.Range("K86").Value = 3160.6667
$ws.Range("L86").Value = 4840.364
$ws.Range("M86").Value = -2037.6667
$ws.Range("N86").Value = -7086.364
$ws.Range("H89").Value = 4084.5
$ws.Range("I89").Value = 3160.6667
$ws.Range("J89").Value = 4840.364
$ws.Range("K89").Value = 15803.3335
$ws.Range("L89").Value = 24201.82
$ws.Range("M89").Value = -10187.3335
$ws.Range("N89").Value = -35433.82
$ws.Range("H107").Value = 45454996
$ws.Range("I107").Value = 490.26315
$ws.Range("J107").Value = 333333540
$ws.Range("K107").Value = 490.26315
$ws.Range("L107").Value = 333333540
$ws.Range("M107").Value = 1429.73685
$ws.Range("N107").Value = -333337380
$ws.Range("H116").Value = 2860
$ws.Range("J116").Value = 3000
$ws.Range("L116").Value = 3000
$ws.Range("N116").Value = -9884
$ws.Range("H131").Value = 90910380
$ws.Range("I131").Value = 100001220
$ws.Range("J131").Value = 2000
$ws.Range("K131").Value = 300003660
$ws.Range("L131").Value = 6000
$ws.Range("M131").Value = -299998620
$ws.Range("N131").Value = -16080
$ws.Range("H132").Value = 4028.7104
$ws.Range("I132").Value = 1442.4642
$ws.Range("J132").Value = 11270.2
$ws.Range("K132").Value = 4327.392599999999
$ws.Range("L132").Value = 33810.60000000001
$ws.Range("M132").Value = -1797.392599999999
$ws.Range("N132").Value = -38870.60000000001
$ws.Range("H138").Value = 2848.4167
$ws.Range("I138").Value = 2355.6775
$ws.Range("J138").Value = 5903.4
$ws.Range("K138").Value = 7067.032499999999
$ws.Range("L138").Value = 17710.2
$ws.Range("M138").Value = -1927.032499999999
$ws.Range("N138").Value = -27990.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30560.242
$ws.Range("I32").Value = 14380
$ws.Range("J32").Value = 47751.75
$ws.Range("K32").Value = 14380
$ws.Range("L32").Value = 47751.75
$ws.Range("M32").Value = -14093
$ws.Range("N32").Value = -48325.75
$ws.Range("H45").Value = 1268.875
$ws.Range("I45").Value = 1230.9231
$ws.Range("J45").Value = 1433.3334
$ws.Range("K45").Value = 1230.9231
$ws.Range("L45").Value = 1433.3334
$ws.Range("M45").Value = -853.9231
$ws.Range("N45").Value = -2187.3334
$ws.Range("H97").Value = 1067.8948
$ws.Range("I97").Value = 972.6667
$ws.Range("J97").Value = 1425
$ws.Range("K97").Value = 972.6667
$ws.Range("L97").Value = 1425
$ws.Range("M97").Value = -476.6667
$ws.Range("N97").Value = -2417
$ws.Range("H127").Value = 34698.57
$ws.Range("J127").Value = 34698.57
$ws.Range("L127").Value = 34698.57
$ws.Range("N127").Value = -44618.57

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 24299.666
$ws.Range("J40").Value = 24299.666
$ws.Range("L40").Value = 24299.666
$ws.Range("N40").Value = -24829.666
$ws.Range("H86").Value = 1854.85
$ws.Range("I86").Value = 1775.3846
$ws.Range("J86").Value = 2002.4286
$ws.Range("K86").Value = 1775.3846
$ws.Range("L86").Value = 2002.4286
$ws.Range("M86").Value = -652.3846000000001
$ws.Range("N86").Value = -4248.4286
$ws.Range("H89").Value = 1854.85
$ws.Range("I89").Value = 1775.3846
$ws.Range("J89").Value = 2002.4286
$ws.Range("K89").Value = 8876.923000000001
$ws.Range("L89").Value = 10012.143
$ws.Range("M89").Value = -3260.923000000001
$ws.Range("N89").Value = -21244.143
$ws.Range("H127").Value = 30000
$ws.Range("J127").Value = 30000
$ws.Range("L127").Value = 30000
$ws.Range("N127").Value = -39920
$ws.Range("H134").Value = 916804.3
$ws.Range("I134").Value = 1371611.5
$ws.Range("J134").Value = 7189.8667
$ws.Range("K134").Value = 4114834.5
$ws.Range("L134").Value = 21569.6001
$ws.Range("M134").Value = -4112299.5
$ws.Range("N134").Value = -26639.6001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3773.7
$ws.Range("I31").Value = 2315.1538
$ws.Range("J31").Value = 4889.0586
$ws.Range("K31").Value = 2315.1538
$ws.Range("L31").Value = 4889.0586
$ws.Range("M31").Value = -2020.1538
$ws.Range("N31").Value = -5479.0586
$ws.Range("H34").Value = 3773.7
$ws.Range("I34").Value = 2315.1538
$ws.Range("J34").Value = 4889.0586
$ws.Range("K34").Value = 2315.1538
$ws.Range("L34").Value = 4889.0586
$ws.Range("M34").Value = -2113.1538
$ws.Range("N34").Value = -5293.0586
$ws.Range("H62").Value = 8180
$ws.Range("I62").Value = 9475
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 9475
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -8851
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 8180
$ws.Range("I65").Value = 9475
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 47375
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -44255
$ws.Range("N65").Value = -21240
$ws.Range("H122").Value = 1250
$ws.Range("I122").Value = 976
$ws.Range("J122").Value = 1592.5
$ws.Range("K122").Value = 2928
$ws.Range("L122").Value = 4777.5
$ws.Range("M122").Value = -478
$ws.Range("N122").Value = -9677.5
$ws.Range("H141").Value = 74518.25999999999
$ws.Range("I141").Value = 19950
$ws.Range("J141").Value = 79715.234
$ws.Range("K141").Value = 19950
$ws.Range("L141").Value = 79715.234
$ws.Range("M141").Value = -14770
$ws.Range("N141").Value = -90075.234

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2466.6667
$ws.Range("I80").Value = 2450
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 2450
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -1452
$ws.Range("N80").Value = -4496
$ws.Range("H83").Value = 2466.6667
$ws.Range("I83").Value = 2450
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 12250
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -7258
$ws.Range("N83").Value = -22484
$ws.Range("H126").Value = 1576
$ws.Range("I126").Value = 1270.1
$ws.Range("J126").Value = 1867.3334
$ws.Range("K126").Value = 3810.3
$ws.Range("L126").Value = 5602.0002
$ws.Range("M126").Value = -1340.3
$ws.Range("N126").Value = -10542.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 687.625
$ws.Range("I16").Value = 645.0909
$ws.Range("J16").Value = 1155.5
$ws.Range("K16").Value = 645.0909
$ws.Range("L16").Value = 1155.5
$ws.Range("M16").Value = -475.0909
$ws.Range("N16").Value = -1495.5
$ws.Range("H61").Value = 3666.6667
$ws.Range("I61").Value = 3000
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 3000
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -2798
$ws.Range("N61").Value = -4404
$ws.Range("H68").Value = 19627.455
$ws.Range("I68").Value = 21290.2
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 21290.2
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -20541.2
$ws.Range("N68").Value = -4498
$ws.Range("H71").Value = 19627.455
$ws.Range("I71").Value = 21290.2
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 106451
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -102707
$ws.Range("N71").Value = -22488
$ws.Range("H100").Value = 2042.6
$ws.Range("I100").Value = 1740.8
$ws.Range("J100").Value = 2344.4
$ws.Range("K100").Value = 1740.8
$ws.Range("L100").Value = 2344.4
$ws.Range("M100").Value = -1199.8
$ws.Range("N100").Value = -3426.4
$ws.Range("H113").Value = 3666.6667
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -8340
$ws.Range("H122").Value = 3099.9443
$ws.Range("I122").Value = 2941.5833
$ws.Range("J122").Value = 3416.6667
$ws.Range("K122").Value = 8824.749899999999
$ws.Range("L122").Value = 10250.0001
$ws.Range("M122").Value = -6374.749899999999
$ws.Range("N122").Value = -15150.0001
$ws.Range("H136").Value = 7785.15
$ws.Range("I136").Value = 1990
$ws.Range("J136").Value = 13580.3
$ws.Range("K136").Value = 5970
$ws.Range("L136").Value = 40740.89999999999
$ws.Range("M136").Value = -3420
$ws.Range("N136").Value = -45840.89999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1009.3077
$ws.Range("I100").Value = 568
$ws.Range("J100").Value = 2002.25
$ws.Range("K100").Value = 1136
$ws.Range("L100").Value = 4004.5
$ws.Range("M100").Value = -595
$ws.Range("N100").Value = -5086.5
$ws.Range("H107").Value = 3039.8718
$ws.Range("I107").Value = 3857.0688
$ws.Range("J107").Value = 670
$ws.Range("K107").Value = 11571.2064
$ws.Range("L107").Value = 2010
$ws.Range("M107").Value = -9651.206399999999
$ws.Range("N107").Value = -5850
$ws.Range("H122").Value = 42627.195
$ws.Range("I122").Value = 5298.9
$ws.Range("J122").Value = 110496.82
$ws.Range("K122").Value = 15896.7
$ws.Range("L122").Value = 331490.46
$ws.Range("M122").Value = -13446.7
$ws.Range("N122").Value = -336390.46
$ws.Range("H126").Value = 931.6667
$ws.Range("I126").Value = 930.2222
$ws.Range("J126").Value = 936
$ws.Range("K126").Value = 2790.6666
$ws.Range("L126").Value = 2808
$ws.Range("M126").Value = -320.6666
$ws.Range("N126").Value = -7748
$ws.Range("H127").Value = 20000
$ws.Range("J127").Value = 20000
$ws.Range("L127").Value = 20000
$ws.Range("N127").Value = -29920
